$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5599.727
$ws.Range("I40").Value = 6311
$ws.Range("J40").Value = 2399
$ws.Range("K40").Value = 6311
$ws.Range("L40").Value = 2399
$ws.Range("M40").Value = -6136
$ws.Range("N40").Value = -2749
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").ClearContents()
$ws.Range("H112").Value = 1616.6897
$ws.Range("J112").Value = 1648.509
$ws.Range("L112").Value = 4945.527
$ws.Range("N112").Value = -7161.527
$ws.Range("H132").Value = 2003.7963
$ws.Range("I132").Value = 2005.88
$ws.Range("K132").Value = 6017.64
$ws.Range("M132").Value = -3487.64
$ws.Range("H135").Value = 892.1818
$ws.Range("J135").Value = 1975
$ws.Range("L135").Value = 17775
$ws.Range("N135").Value = -22845
$ws.Range("H137").Value = 2015.8
$ws.Range("I137").Value = 1134.3954
$ws.Range("J137").Value = 3738.5454
$ws.Range("K137").Value = 3403.1862
$ws.Range("L137").Value = 11215.6362
$ws.Range("M137").Value = -853.1862000000001
$ws.Range("N137").Value = -16315.6362
$ws.Range("H138").Value = 4175.467
$ws.Range("J138").Value = 5105
$ws.Range("L138").Value = 15315
$ws.Range("N138").Value = -25595
$ws.Range("H141").Value = 1462.4783
$ws.Range("I141").Value = 1428.9546
$ws.Range("K141").Value = 4286.8638
$ws.Range("M141").Value = 893.1361999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1095.2388
$ws.Range("I32").Value = 879.069
$ws.Range("K32").Value = 879.069
$ws.Range("M32").Value = -592.069
$ws.Range("H61").Value = 13488.272
$ws.Range("I61").Value = 22477.6
$ws.Range("K61").Value = 22477.6
$ws.Range("M61").Value = -22265.6
$ws.Range("H111").Value = 100644
$ws.Range("J111").Value = 100644
$ws.Range("L111").Value = 100644
$ws.Range("N111").Value = -108824
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H136").Value = 13488.272
$ws.Range("I136").Value = 22477.6
$ws.Range("K136").Value = 67432.79999999999
$ws.Range("M136").Value = -64882.79999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1814.0625
$ws.Range("I20").Value = 1802
$ws.Range("J20").Value = 1866.3334
$ws.Range("K20").Value = 1802
$ws.Range("L20").Value = 1866.3334
$ws.Range("M20").Value = -1555
$ws.Range("N20").Value = -2360.3334
$ws.Range("H86").Value = 2310.6667
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 2310.6667
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H105").Value = 2766.375
$ws.Range("I105").Value = 2447.4285
$ws.Range("K105").Value = 2447.4285
$ws.Range("M105").Value = -700.4285
$ws.Range("H134").Value = 2595.25
$ws.Range("I134").Value = 1783.3334
$ws.Range("K134").Value = 5350.0002
$ws.Range("M134").Value = -2815.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 643.625
$ws.Range("I22").Value = 649.86957
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 649.86957
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -299.86957
$ws.Range("N22").Value = -1200
$ws.Range("H31").Value = 4963.684
$ws.Range("I31").Value = 9073
$ws.Range("J31").Value = 3067.077
$ws.Range("K31").Value = 9073
$ws.Range("L31").Value = 3067.077
$ws.Range("M31").Value = -8778
$ws.Range("N31").Value = -3657.077
$ws.Range("H34").Value = 4963.684
$ws.Range("I34").Value = 9073
$ws.Range("J34").Value = 3067.077
$ws.Range("K34").Value = 9073
$ws.Range("L34").Value = 3067.077
$ws.Range("M34").Value = -8871
$ws.Range("N34").Value = -3471.077
$ws.Range("H58").Value = 3580
$ws.Range("I58").Value = 1937.2858
$ws.Range("K58").Value = 1937.2858
$ws.Range("M58").Value = -1734.2858
$ws.Range("H86").Value = 7482.6665
$ws.Range("I86").Value = 5500
$ws.Range("J86").Value = 8474
$ws.Range("K86").Value = 5500
$ws.Range("L86").Value = 8474
$ws.Range("M86").Value = -4377
$ws.Range("N86").Value = -10720
$ws.Range("H89").Value = 7482.6665
$ws.Range("I89").Value = 5500
$ws.Range("J89").Value = 8474
$ws.Range("K89").Value = 27500
$ws.Range("L89").Value = 42370
$ws.Range("M89").Value = -21884
$ws.Range("N89").Value = -53602
$ws.Range("H99").Value = 2224.375
$ws.Range("I99").Value = 1799.5
$ws.Range("K99").Value = 1799.5
$ws.Range("M99").Value = -301.5
$ws.Range("H105").Value = 6548.3335
$ws.Range("I105").Value = 7822.75
$ws.Range("K105").Value = 7822.75
$ws.Range("M105").Value = -6075.75
$ws.Range("H122").Value = 2298.2
$ws.Range("J122").Value = 3151.889
$ws.Range("L122").Value = 9455.667000000001
$ws.Range("N122").Value = -14355.667
$ws.Range("H126").Value = 2224.375
$ws.Range("I126").Value = 1799.5
$ws.Range("K126").Value = 5398.5
$ws.Range("M126").Value = -2928.5
$ws.Range("H132").Value = 4228.0835
$ws.Range("I132").Value = 2527.875
$ws.Range("J132").Value = 7628.5
$ws.Range("K132").Value = 7583.625
$ws.Range("L132").Value = 22885.5
$ws.Range("M132").Value = -5053.625
$ws.Range("N132").Value = -27945.5
$ws.Range("H134").Value = 1586.0667
$ws.Range("I134").Value = 1585.1428
$ws.Range("K134").Value = 4755.428400000001
$ws.Range("M134").Value = -2220.428400000001
$ws.Range("H136").Value = 3580
$ws.Range("I136").Value = 1937.2858
$ws.Range("K136").Value = 5811.857400000001
$ws.Range("M136").Value = -3261.857400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 1675337.2
$ws.Range("I8").Value = 1675337.2
$ws.Range("K8").Value = 5026011.6
$ws.Range("M8").Value = -5025872.6
$ws.Range("H92").Value = 396.7
$ws.Range("I92").Value = 207.25
$ws.Range("K92").Value = 621.75
$ws.Range("M92").Value = 626.25
$ws.Range("H102").Value = 9144.4
$ws.Range("J102").Value = 10393.385
$ws.Range("L102").Value = 31180.155
$ws.Range("N102").Value = -36048.155
$ws.Range("H121").Value = 803.125
$ws.Range("J121").Value = 862.7273
$ws.Range("L121").Value = 2588.1819
$ws.Range("N121").Value = -5208.1819
$ws.Range("H131").Value = 1289.7894
$ws.Range("J131").Value = 2006.4
$ws.Range("L131").Value = 6019.200000000001
$ws.Range("N131").Value = -16099.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 47777.5
$ws.Range("J32").Value = 47777.5
$ws.Range("L32").Value = 47777.5
$ws.Range("N32").Value = -48369.5
$ws.Range("H45").Value = 26517.857
$ws.Range("I45").Value = 10000
$ws.Range("J45").Value = 29270.834
$ws.Range("K45").Value = 10000
$ws.Range("L45").Value = 29270.834
$ws.Range("M45").Value = -9441
$ws.Range("N45").Value = -30388.834
$ws.Range("H113").Value = 5444.4443
$ws.Range("I113").Value = 5708.3335
$ws.Range("K113").Value = 5708.3335
$ws.Range("M113").Value = -3538.3335
$ws.Range("H122").Value = 5611.7144
$ws.Range("I122").Value = 7253.357
$ws.Range("K122").Value = 21760.071
$ws.Range("M122").Value = -19310.071
$ws.Range("H132").Value = 29208.584
$ws.Range("I132").Value = 15610.223
$ws.Range("J132").Value = 70003.664
$ws.Range("K132").Value = 46830.669
$ws.Range("L132").Value = 210010.992
$ws.Range("M132").Value = -44300.669
$ws.Range("N132").Value = -215070.992

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 1258.3334
$ws.Range("I12").Value = 1000
$ws.Range("J12").Value = 1775
$ws.Range("K12").Value = 1000
$ws.Range("L12").Value = 1775
$ws.Range("M12").Value = -830
$ws.Range("N12").Value = -2115
$ws.Range("H16").Value = 15549.875
$ws.Range("I16").Value = 22580
$ws.Range("J16").Value = 3833
$ws.Range("K16").Value = 22580
$ws.Range("L16").Value = 3833
$ws.Range("M16").Value = -22410
$ws.Range("N16").Value = -4173
$ws.Range("H55").Value = 314
$ws.Range("I55").Value = 315
$ws.Range("J55").Value = 313.33334
$ws.Range("K55").Value = 315
$ws.Range("L55").Value = 313.33334
$ws.Range("M55").Value = -142
$ws.Range("N55").Value = -659.33334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 60078.75
$ws.Range("J70").Value = 60078.75
$ws.Range("L70").Value = 60078.75
$ws.Range("N70").Value = -60708.75
$ws.Range("H73").Value = 60078.75
$ws.Range("J73").Value = 60078.75
$ws.Range("L73").Value = 60078.75
$ws.Range("N73").Value = -62262.75
$ws.Range("H122").Value = 1482.841
$ws.Range("I122").Value = 1376.8462
$ws.Range("K122").Value = 4130.5386
$ws.Range("M122").Value = -1680.5386
$ws.Range("H132").Value = 3301.2983
$ws.Range("I132").Value = 1977.9773
$ws.Range("J132").Value = 7780.231
$ws.Range("K132").Value = 5933.9319
$ws.Range("L132").Value = 23340.693
$ws.Range("M132").Value = -3403.9319
$ws.Range("N132").Value = -28400.693
$ws.Range("H136").Value = 2915.8215
$ws.Range("I136").Value = 2325.476
$ws.Range("J136").Value = 4686.857
$ws.Range("K136").Value = 6976.428
$ws.Range("L136").Value = 14060.571
$ws.Range("M136").Value = -4426.428
$ws.Range("N136").Value = -19160.571
